# Revise abc.xlsx: add a new row (row 4) on Sheet1 with a "now " label,
# a date value (2018-07-31) and a time value (15:04), matching the
# "20180731 15:04 revised abc.xlsx" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in E4 (second shared string "now ").
$ws.Range("E4").Value = "now "

# G4 first, so its new style (numFmtId 20 / h:mm) becomes cellXfs index 1.
$ws.Range("G4").NumberFormat = "h:mm"
$ws.Range("G4").Value = 0.62777777777777777

# F4 second, so its new style (numFmtId 14 / mm-dd-yy) becomes cellXfs index 2.
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("F4").Value = 43312

# Widen column F to fit the date value (bestFit-style width).
$ws.Columns.Item(6).ColumnWidth = 9.5

# Select G4, matching the saved selection in the workbook.
$null = $ws.Range("G4").Select()
